# Updates crypto price/volume data in the "cryptos" worksheet to the
# latest scraped values (per GitHub Actions commit "Updated cryptos list").
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
#
# For numeric-looking Price values (single decimal point, e.g. "6.10"),
# we force the cell to store a literal text string (not a number) by
# temporarily applying a Text number format, then restoring the cell's
# style to "Normal" so no residual formatting/style difference remains.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '61.509.06'
$ws.Cells.Item(2, 5).Value = '  +0.86%  '

$ws.Cells.Item(3, 4).Value = '3.385.60'
$ws.Cells.Item(3, 5).Value = '  -0.35%  '

$ws.Cells.Item(4, 5).Value = '  -0.03%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '574.87'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.33%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '140.56'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -1.52%  '

$ws.Cells.Item(8, 5).Value = '  -0.61%  '

$ws.Cells.Item(9, 5).Value = '  +1.68%  '

$ws.Cells.Item(10, 5).Value = '  -1.48%  '

$ws.Cells.Item(11, 5).Value = '  -2.77%  '

$ws.Cells.Item(12, 4).Value = '3.964.81'
$ws.Cells.Item(12, 5).Value = '  -0.28%  '

$ws.Cells.Item(13, 5).Value = '  +0.30%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '28.37'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +1.07%  '

$ws.Cells.Item(15, 4).Value = '3.385.93'
$ws.Cells.Item(15, 5).Value = '  -0.22%  '

$ws.Cells.Item(16, 5).Value = '  -0.92%  '

$ws.Cells.Item(17, 4).Value = '61.524.22'
$ws.Cells.Item(17, 5).Value = '  +0.76%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '6.10'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -0.71%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '13.58'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -2.07%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '8.99'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.17%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '390.53'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +1.70%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '74.97'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.96%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.551'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -1.49%  '

$ws.Cells.Item(24, 5).Value = '  -0.04%  '

$ws.Cells.Item(25, 5).Value = '  -5.07%  '

$ws.Cells.Item(26, 5).Value = '  +7.17%  '

$ws.Cells.Item(27, 5).Value = '  +0.03%  '

$ws.Cells.Item(28, 5).Value = '  -2.22%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '8.03'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +0.14%  '

$ws.Cells.Item(30, 5).Value = '  -1.04%  '

$ws.Cells.Item(31, 5).Value = '  +0.00%  '

$ws.Cells.Item(32, 5).Value = '  -1.49%  '

$ws.Cells.Item(33, 5).Value = '  -1.40%  '

$ws.Cells.Item(34, 2).Value = 'Monero'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '169.23'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +1.00%  '

$ws.Cells.Item(35, 2).Value = 'Aptos'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '6.89'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -1.71%  '

$ws.Cells.Item(36, 5).Value = '  +0.32%  '

$ws.Cells.Item(37, 4).Value = '3.419.51'
$ws.Cells.Item(37, 5).Value = '  -0.24%  '

$ws.Cells.Item(38, 5).Value = '  -1.83%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.0765'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.94%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '26.17'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -5.19%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.779'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -0.37%  '

$ws.Cells.Item(42, 5).Value = '  -0.64%  '

$ws.Cells.Item(43, 5).Value = '  -1.60%  '

$ws.Cells.Item(44, 5).Value = '  +1.39%  '

$ws.Cells.Item(45, 4).Value = '2.459.37'
$ws.Cells.Item(45, 5).Value = '  -1.04%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '22.85'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -0.65%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '6.64'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -2.55%  '

$ws.Cells.Item(48, 5).Value = '  +0.02%  '

$ws.Cells.Item(49, 5).Value = '  -1.13%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '2.00'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -5.50%  '

$ws.Cells.Item(51, 5).Value = '  -2.22%  '
